$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Change 1: italicize the file names in the first paragraph, keeping
# " y " between the two ".docx" filenames as plain text.
# ------------------------------------------------------------------

# "[ingeS]SpmpVersion2.1" -> italic
$r = $d.Content
$null = $r.Find.Execute("[ingeS]SpmpVersion2.1", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r.Font.Italic = $true

# "(Linea_base)" -> italic
$r = $d.Content
$null = $r.Find.Execute("(Linea_base)", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r.Font.Italic = $true

# ".docx y CasosdeUso1.0.docx" -> split into ".docx" (italic) + " y " (plain)
# + "CasosdeUso1.0.docx" (italic), without ever touching the " y " piece so
# it keeps no direct formatting at all (matching the target markup).
$r = $d.Content
$null = $r.Find.Execute(".docx y CasosdeUso1.0.docx", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$start = $r.Start
$end = $r.End

$firstDocx = $d.Range($start, $start + 5)
$firstDocx.Font.Italic = $true

$secondDocx = $d.Range($end - 18, $end)
$secondDocx.Font.Italic = $true

# ------------------------------------------------------------------
# Change 2: split out "Súper" (correcting "Super") into its own run,
# between "ción del proyecto " and " Triumph ".
# ------------------------------------------------------------------

# Fix the typo first. This is a text-content edit, so the engine
# re-coalesces every same-formatted run across the whole paragraph
# (collapsing "ción del proyecto ", "Súper" and " Triumph " together,
# right up against the following "realizado por IMind" run).
$r = $d.Content
$null = $r.Find.Execute("Super", $true, $false, $false, $false, $false, $true, 1, $false, "Súper", 2)

# Now re-establish the run boundaries the diff expects, purely through
# (idempotent) direct-formatting toggles, which only affect the local
# run boundary instead of re-coalescing the whole paragraph again.

# Boundary between "...planea" and "ción del proyecto ":
$r = $d.Content
$null = $r.Find.Execute("ción del proyecto ", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r.Font.Bold = $true
$r.Font.Bold = $false

# Boundary around "Súper" itself:
$r = $d.Content
$null = $r.Find.Execute("Súper", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r.Font.Bold = $true
$r.Font.Bold = $false

# Boundary between " Triumph " and "realizado por IMind":
$r = $d.Content
$null = $r.Find.Execute(" Triumph ", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r.Font.Bold = $true
$r.Font.Bold = $false

# ------------------------------------------------------------------
# Change 3: remove the trailing "Falta el plan de mejoras..." text
# (and the stray single-space run before it).
# ------------------------------------------------------------------
$r = $d.Content
$null = $r.Find.Execute(" Falta el plan de mejoras del proceso, que está incompleto pero preferiblemente se colocará en la próxima versión complemente. ", $false, $false, $false, $false, $false, $true, 1, $false, "", 2)
